$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the "Team Name" column (column B) ---
# Clear the header + data that used to live in column B.
$ws.Range("B1:B5").ClearContents()

# Column B no longer needs the box-drawing border it had as the right-hand
# wall of the old two-column table; strip borders on all sides for B1:B5.
for ($r = 1; $r -le 5; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $cell.Borders.Item(7).LineStyle  = -4142   # left
    $cell.Borders.Item(8).LineStyle  = -4142   # top
    $cell.Borders.Item(9).LineStyle  = -4142   # bottom
    $cell.Borders.Item(10).LineStyle = -4142   # right
}

# Column A now stands alone, so its box needs a right edge that used to be
# supplied implicitly by column B's left edge. Rebuild the box using only
# column A: A1 keeps a full box, A2/A3 get left+right, A4 closes the box
# with left+right+bottom.
$a1 = $ws.Range("A1")
$a1.BorderAround(1, 2)

$a2 = $ws.Range("A2")
$a2.BorderAround(1, 2)
$a2.Borders.Item(8).LineStyle  = -4142   # no top (A1's bottom edge closes it)
$a2.Borders.Item(9).LineStyle  = -4142   # no bottom

$a3 = $ws.Range("A3")
$a3.BorderAround(1, 2)
$a3.Borders.Item(8).LineStyle  = -4142   # no top
$a3.Borders.Item(9).LineStyle  = -4142   # no bottom

$a4 = $ws.Range("A4")
$a4.BorderAround(1, 2)
$a4.Borders.Item(8).LineStyle  = -4142   # no top

# Column B keeps its width, but no longer carries any special column style.
$ws.Columns("B").ColumnWidth = 25.5703125

# Move the active selection the way the source workbook shows it post-edit.
$ws.Range("B10").Select()
